$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 369 (existing rows 369:434 shift down to 371:436)
$ws.Range("A369:A370").EntireRow.Insert()

# New row 369: "Primera" quality record for date 44504
$ws.Range("A369").Value = 8
$ws.Range("B369").Value = "Terminal La Palmera de La Serena"
$ws.Range("C369").Value = "Coquimbo"
$ws.Range("D369").Value = 44504
$ws.Range("E369").Value = 4
$ws.Range("F369").Value = 100112008
$ws.Range("G369").Value = "Coliflor"
$ws.Range("H369").Value = "Sin especificar"
$ws.Range("I369").Value = "Primera"
$ws.Range("J369").Value = 2200
$ws.Range("K369").Value = 600
$ws.Range("L369").Value = 700
$ws.Range("M369").Value = 650
$ws.Range("N369").Value = "`$/unidad"
$ws.Range("O369").Value = "Provincia del Elquí"
$ws.Range("P369").Value = 650
$ws.Range("Q369").Value = 1
$ws.Range("R369").Value = "Hortaliza"

# New row 370: "Segunda" quality record for date 44504
$ws.Range("A370").Value = 8
$ws.Range("B370").Value = "Terminal La Palmera de La Serena"
$ws.Range("C370").Value = "Coquimbo"
$ws.Range("D370").Value = 44504
$ws.Range("E370").Value = 4
$ws.Range("F370").Value = 100112008
$ws.Range("G370").Value = "Coliflor"
$ws.Range("H370").Value = "Sin especificar"
$ws.Range("I370").Value = "Segunda"
$ws.Range("J370").Value = 1320
$ws.Range("K370").Value = 500
$ws.Range("L370").Value = 550
$ws.Range("M370").Value = 525
$ws.Range("N370").Value = "`$/unidad"
$ws.Range("O370").Value = "Provincia del Elquí"
$ws.Range("P370").Value = 525
$ws.Range("Q370").Value = 1
$ws.Range("R370").Value = "Hortaliza"

# Match the date cell number format used by the rest of column D
$ws.Range("D369:D370").NumberFormat = $ws.Range("D371").NumberFormat
